$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.516.28"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "3.493.10"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'613.69"
$ws.Range("E5").Value = "  +5.56%  "
$ws.Range("D6").Value = "'189.52"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -4.65%  "
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("E11").Value = "  -3.19%  "
$ws.Range("E12").Value = "  -3.90%  "
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "4.057.72"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").Value = "'599.94"
$ws.Range("E15").Value = "  +4.08%  "
$ws.Range("D16").Value = "69.584.22"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("D17").Value = "'18.99"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "'12.58"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("D19").Value = "3.482.38"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").Value = "'0.986"
$ws.Range("E21").Value = "  -1.76%  "
$ws.Range("D22").Value = "'17.12"
$ws.Range("E22").Value = "  -3.05%  "
$ws.Range("D23").Value = "'105.66"
$ws.Range("E23").Value = "  +12.59%  "
$ws.Range("D24").Value = "'4.71"
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("D25").Value = "'5.12"
$ws.Range("E25").Value = "  +5.01%  "
$ws.Range("E26").Value = "  +2.64%  "
$ws.Range("E27").Value = "  -2.67%  "
$ws.Range("D28").Value = "'9.69"
$ws.Range("E28").Value = "  +4.18%  "
$ws.Range("D29").Value = "'33.39"
$ws.Range("E29").Value = "  +2.50%  "
$ws.Range("E30").Value = "  -3.36%  "
$ws.Range("D31").Value = "'4.17"
$ws.Range("E31").Value = "  +10.07%  "
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").Value = "'63.34"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("E35").Value = "  -5.28%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "'3.66"
$ws.Range("E37").Value = "  +6.07%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "3.622.72"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("E39").Value = "  -4.19%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'36.74"
$ws.Range("E40").Value = "  -3.87%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'503.97"
$ws.Range("E41").Value = "  -6.69%  "
$ws.Range("D42").Value = "0.0₃0773"
$ws.Range("E42").Value = "  -4.49%  "
$ws.Range("D43").Value = "'0.137"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").Value = "'3.32"
$ws.Range("E47").Value = "  -4.96%  "
$ws.Range("D48").Value = "'8.76"
$ws.Range("E48").Value = "  -6.20%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").Value = "'131.63"
$ws.Range("E50").Value = "  -2.38%  "
